$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# "Latest HO Xliff Generate Date" on Overview sheet (shared with de-de's
# "Correspond Handoff Datetime" column, same underlying string).
$wsOverview.Range("G2").Value = "2016-09-04 06:19:35"
$wsOverview.Range("G4").Value = "2016-09-04 06:19:35"

# de-de sheet: same timestamp string referenced via "Correspond Handoff Datetime".
$wsDeDe.Range("H2").Value = "2016-09-04 06:19:35"
$wsDeDe.Range("H4").Value = "2016-09-04 06:19:35"

# "Priority" column ("ht" -> "mt") on both zh-cn and de-de sheets.
$wsZhCn.Range("E2").Value = "mt"
$wsZhCn.Range("E4").Value = "mt"
$wsDeDe.Range("E2").Value = "mt"
$wsDeDe.Range("E4").Value = "mt"

# zh-cn sheet "Correspond Handoff Datetime" column.
$wsZhCn.Range("H2").Value = "2016-09-04 06:19:31"
$wsZhCn.Range("H4").Value = "2016-09-04 06:19:31"

# zh-cn sheet "Correspond Handback DateTime" column.
$wsZhCn.Range("K2").Value = "2016-09-04 06:19:47"
$wsZhCn.Range("K4").Value = "2016-09-04 06:19:47"

# de-de sheet "Correspond Handback DateTime" column.
$wsDeDe.Range("K2").Value = "2016-09-04 06:19:54"
$wsDeDe.Range("K4").Value = "2016-09-04 06:19:54"
